$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Style the header row (A1:F1): bold font, thin border all around,
#    centered horizontally, top-aligned vertically.
#    Configure a single cell fully, then propagate the resulting format
#    (not the values) to the rest of the header row.
# ---------------------------------------------------------------------------
$headerCell = $ws.Cells.Item(1, 1)
$headerCell.Font.Bold = $true
$headerCell.Borders.LineStyle = 1
$headerCell.HorizontalAlignment = -4108   # xlCenter
$headerCell.VerticalAlignment = -4160     # xlTop

$headerCell.Copy()
$ws.Range("B1:F1").PasteSpecial(-4122)    # xlPasteFormats

# ---------------------------------------------------------------------------
# 2) Remove the old job rows 13-21 (they're replaced / no longer listed).
# ---------------------------------------------------------------------------
$ws.Range("A13:F21").EntireRow.Delete()

# ---------------------------------------------------------------------------
# 3) Update row 12 with the new job listing.
# ---------------------------------------------------------------------------
$ws.Range("A12").Value = "Golang Architect / Principal Backend Architect - Atlanta, GA"
$ws.Range("B12").Value = "https://www.dice.com/job-detail/d909bde8-866c-4ca7-8874-8e2c93f3aad7"
$ws.Range("C12").Value = "Atlanta, Georgia"
$ws.Range("D12").Value = "Third Party, Contract"
$ws.Range("E12").Value = '$58 - $68 per hour'
$ws.Range("F12").Value = "Bayside Solutions"
